$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (2021) is new. Give A5 the same header-year style as A2:A4 (bold, bordered, centered)
# by copying row 4's formatting before writing the new values.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "2021年"
$ws.Range("B5").Value = 9.1
$ws.Range("C5").Value = 29.2
$ws.Range("D5").Value = -20.7
$ws.Range("E5").Value = 72.5
$ws.Range("F5").Value = -0.8
$ws.Range("G5").Value = 67.2
$ws.Range("H5").Value = 8.300000000000001
$ws.Range("I5").Value = 32.1
$ws.Range("J5").Value = 31
$ws.Range("K5").Value = -9.300000000000001
$ws.Range("L5").Value = -10.6
$ws.Range("M5").Value = 3.4
$ws.Range("N5").Value = -39.8
$ws.Range("O5").Value = -1
$ws.Range("P5").Value = 12.6
$ws.Range("Q5").Value = -6.4
$ws.Range("R5").Value = -46.6
$ws.Range("S5").Value = 79.8
$ws.Range("T5").Value = 8.199999999999999
$ws.Range("U5").Value = 12
$ws.Range("V5").Value = 16.2
$ws.Range("W5").Value = 37.7
$ws.Range("X5").Value = 10.2
$ws.Range("Y5").Value = 1.1
$ws.Range("Z5").Value = 0.6
$ws.Range("AA5").Value = 26.9
$ws.Range("AB5").Value = 29.6
$ws.Range("AC5").Value = 31.7
$ws.Range("AD5").Value = -5.5
$ws.Range("AE5").Value = 14.8
$ws.Range("AF5").Value = -42
$ws.Range("AG5").Value = "'"
$ws.Range("AG5").ClearFormats() | Out-Null
$ws.Range("AH5").Value = -18.7
$ws.Range("AI5").Value = -18
$ws.Range("AJ5").Value = -14.5
$ws.Range("AK5").Value = 8.800000000000001
$ws.Range("AL5").Value = 7.5
$ws.Range("AM5").Value = -13.7
$ws.Range("AN5").Value = -7
$ws.Range("AO5").Value = -7.6
$ws.Range("AP5").Value = -13.5
$ws.Range("AQ5").Value = 0.1
$ws.Range("AR5").Value = -14.7
$ws.Range("AS5").Value = 177.5
$ws.Range("AT5").Value = 72.2
$ws.Range("AU5").Value = 55.9
$ws.Range("AV5").Value = 10.8
$ws.Range("AW5").Value = -17.3
$ws.Range("AX5").Value = -9.9
$ws.Range("AY5").Value = -15.1
$ws.Range("AZ5").Value = 7.3
$ws.Range("BA5").Value = 3.5
$ws.Range("BB5").Value = 10.7
$ws.Range("BC5").Value = -0.5
$ws.Range("BD5").Value = 60.7
$ws.Range("BE5").Value = -22.6
$ws.Range("BF5").Value = -17.8
$ws.Range("BG5").Value = 53.3
$ws.Range("BH5").Value = 14.1
$ws.Range("BI5").Value = "'"
$ws.Range("BI5").ClearFormats() | Out-Null
$ws.Range("BJ5").Value = -14.4
$ws.Range("BK5").Value = 3.6
$ws.Range("BL5").Value = -1.1
$ws.Range("BM5").Value = 0.9
$ws.Range("BN5").Value = -5.3
$ws.Range("BO5").Value = -11.7
$ws.Range("BP5").Value = 53.6
$ws.Range("BQ5").Value = 33.4
$ws.Range("BR5").Value = 10.7
$ws.Range("BS5").Value = 9.300000000000001
$ws.Range("BT5").Value = -4
$ws.Range("BU5").Value = -17.7
$ws.Range("BV5").Value = 0.3
$ws.Range("BW5").Value = 2
$ws.Range("BX5").Value = 30.6
$ws.Range("BY5").Value = 20.6
$ws.Range("BZ5").Value = 10.7
$ws.Range("CA5").Value = 6.1
$ws.Range("CB5").Value = 3.7
$ws.Range("CC5").Value = 10.2
$ws.Range("CD5").Value = -21.1
$ws.Range("CE5").Value = 10.4
$ws.Range("CF5").Value = 17.5
$ws.Range("CG5").Value = 33.3
$ws.Range("CH5").Value = 79.7
$ws.Range("CI5").Value = 17.5
$ws.Range("CJ5").Value = -16
$ws.Range("CK5").Value = 11.2
$ws.Range("CL5").Value = -4.8
$ws.Range("CM5").Value = -4.5
$ws.Range("CN5").Value = 20
$ws.Range("CO5").Value = 34.2
$ws.Range("CP5").Value = 29.8
$ws.Range("CQ5").Value = 13.7
$ws.Range("CR5").Value = -8
$ws.Range("CS5").Value = 43.3
$ws.Range("CT5").Value = -3
$ws.Range("CU5").Value = -9.199999999999999
$ws.Range("CV5").Value = -2.8
$ws.Range("CW5").Value = 36.4
$ws.Range("CX5").Value = 8.5
$ws.Range("CY5").Value = 10.1
$ws.Range("CZ5").Value = 53.7
$ws.Range("DA5").Value = -1.8
$ws.Range("DB5").Value = 14.6
$ws.Range("DC5").Value = 23.1
$ws.Range("DD5").Value = -1
$ws.Range("DE5").Value = -18.3
$ws.Range("DF5").Value = 0.7
$ws.Range("DG5").Value = 80.5
$ws.Range("DH5").Value = 32.7
$ws.Range("DI5").Value = 18.4
$ws.Range("DJ5").Value = -0.3
$ws.Range("DK5").Value = 29.8
